$d = $word.ActiveDocument

function Replace-All($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2) | Out-Null
}

function Replace-InParagraph($paraIndex, $old, $new) {
    $r = $d.Paragraphs($paraIndex).Range
    $r.Find.Execute($old, $true, $false, $false, $false, $false,
                    $true, 1, $false, $new, 2) | Out-Null
}

# Replacing text inside an r:id hyperlink run via Find/Replace (or plain
# Range.Text assignment) drops the run-level color/underline formatting in
# this engine, so restore it explicitly afterwards for such runs.
function Replace-HyperlinkRun($old, $new, $colorBgr, $underline) {
    $r = $d.Content
    $r.Find.Execute($old, $true, $false, $false, $false, $false,
                    $true, 1, $false, "", 0) | Out-Null
    $start = $r.Start
    $end = $r.End
    $target = $d.Range($start, $end)
    $target.Text = $new
    $newEnd = $start + $new.Length
    $target2 = $d.Range($start, $newEnd)
    $target2.Font.Color = $colorBgr
    $target2.Font.Underline = $underline
}

Replace-All "English" "الإنجليزية"
Replace-All " / Portuguese / French / Thai / Vietnamese / Spanish" " /البرتغالية/الفرنسية/التايلندية/الفيتنامية/الإسبانية"
Replace-All "English" "الإنجليزية"
Replace-All "Brief" "المضمون"
Replace-All "An email sent to partners in the target country whose documents failed our verification process. It will be sent via customer.io" "رسالة بريد إلكتروني مرسلة إلى الشركاء في البلد المعني الذين فشلت مستنداتهم في عملية التحقق الخاصة بنا. سيتم إرسالها عبر customer.io"
Replace-All "Target audience" "الجمهور المستهدف"
Replace-All "Invited partners who submitted wrong/incomplete documents" "الشركاء المدعوون الذين قدموا مستندات خاطئة/غير كاملة"
Replace-All "Subject line" "سطر الموضوع"
Replace-All "[EVENT NAME]" "[اسم الحدث]"
Replace-All " — document verification failed " " - فشل التحقق من المستند "
Replace-All "Uh oh! Your documents couldn’t be verified" "أوه! تم التحقق من المستند الخاص بك"
Replace-All "If you have any questions, please contact us via " "إذا كانت لديك أي أسئلة، فاتصل بنا:  "
Replace-All "If you have any questions, please contact your country manager, " "إذا كانت لديك أي أسئلة، فيُرجى الاتصال بمدير بلدك  "
Replace-All "Hi " "مرحبًا  "
Replace-All "We regret to inform you that your documents have failed our verification process as we found the following issues with them: " "نأسف لإبلاغك بأن مستنداتك قد فشلت في عملية التحقق الخاصة بنا حيث وجدنا المشكلات التالية معها: "
Replace-All "A copy of your vaccination certificate" "نسخة من شهادة التطعيم الخاصة بك"
Replace-All ": Document is unclear" ": المستند غير واضح"
Replace-All "Please resubmit the documents above by " "يرجى إعادة تقديم المستندات أعلاه بحلول  "
Replace-All " so we can proceed with the necessary arrangements." " حتى نتمكن من المضي قدمًا في الترتيبات اللازمة."

# "live chat" sits inside an r:id hyperlink run; use the formatting-safe helper
# (1155cc RGB -> BGR 13391121; wdUnderlineSingle = 1)
Replace-HyperlinkRun "live chat" "الدردشة الحية" 13391121 1

# Scoped to paragraph 16 ("Hi [PARTNER NAME], ") to avoid corrupting the
# unrelated ", at " text later in the document
Replace-InParagraph 16 ", " ",، "
